$d = $word.ActiveDocument

function Get-CleanText($para) {
    return $para.Range.Text.TrimEnd([char]13, [char]7)
}

# ------------------------------------------------------------------
# 1) Remove the stray "_GoBack" bookmark from the title paragraph.
#    (It will be recreated, empty, inside the "Goal" heading below.)
# ------------------------------------------------------------------
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

# ------------------------------------------------------------------
# 2) Remove the two empty, italic, indented paragraphs that sit
#    right before the "Goal" heading.
# ------------------------------------------------------------------
$d.Paragraphs.Item(6).Range.Delete()
$d.Paragraphs.Item(6).Range.Delete()

# ------------------------------------------------------------------
# 3) Replace the "Goal" heading paragraph: drop the
#    <w:smartTag element="place"> wrapper around "Goa" (keep the
#    plain text run) and add an empty "_GoBack" bookmark right at
#    the start of the paragraph.
# ------------------------------------------------------------------
$goalPara = $d.Paragraphs.Item(6)
$goalXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:pPr><w:pStyle w:val="Heading3"/></w:pPr>
<w:bookmarkStart w:id="0" w:name="_GoBack"/>
<w:bookmarkEnd w:id="0"/>
<w:r><w:t>Goa</w:t></w:r>
<w:r><w:t>l</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$goalPara.Range.InsertXML($goalXml)

# ------------------------------------------------------------------
# 4) Normalize the <w:attr> order inside the two "date" smartTags
#    (Month, Day, Year instead of Year, Day, Month) in the
#    "June 6 2008 - June 12 2008" paragraph.
# ------------------------------------------------------------------
$dash = [char]0x2013
$dateTarget = " " + $dash + " "
$datePara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ((Get-CleanText $p) -eq $dateTarget) {
        $datePara = $p
        break
    }
}
if ($datePara -eq $null) {
    throw "Could not locate the 'June 6 - June 12' date paragraph"
}
$dateXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:pPr><w:ind w:left="426"/></w:pPr>
<w:smartTag w:uri="urn:schemas-microsoft-com:office:smarttags" w:element="date">
<w:smartTagPr>
<w:attr w:name="Month" w:val="6"/>
<w:attr w:name="Day" w:val="6"/>
<w:attr w:name="Year" w:val="2008"/>
</w:smartTagPr>
<w:r><w:t>June 6</w:t></w:r>
<w:r><w:t xml:space="preserve">, </w:t></w:r>
<w:r><w:t>2008</w:t></w:r>
</w:smartTag>
<w:r><w:t xml:space="preserve"> </w:t></w:r>
<w:r><w:t xml:space="preserve">&#8211; </w:t></w:r>
<w:smartTag w:uri="urn:schemas-microsoft-com:office:smarttags" w:element="date">
<w:smartTagPr>
<w:attr w:name="Month" w:val="6"/>
<w:attr w:name="Day" w:val="12"/>
<w:attr w:name="Year" w:val="2008"/>
</w:smartTagPr>
<w:r><w:t>June 12</w:t></w:r>
<w:r><w:t>, 2008</w:t></w:r>
</w:smartTag>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$datePara.Range.InsertXML($dateXml)

# ------------------------------------------------------------------
# 5) Merge the split "- Label & Goto" / "- Label & Goto in a
#    Diagram" runs (dropping the <w:proofErr> spell-check markers)
#    into single, unbroken runs.
# ------------------------------------------------------------------
$labelPara1 = $null
$labelPara2 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = Get-CleanText $p
    if ($t -eq "- Label & Goto") {
        $labelPara1 = $p
    } elseif ($t -eq "- Label & Goto in a Diagram") {
        $labelPara2 = $p
        break
    }
}
if ($labelPara1 -eq $null -or $labelPara2 -eq $null) {
    throw "Could not locate the 'Label & Goto' paragraphs"
}
$labelRange = $d.Range($labelPara1.Range.Start, $labelPara2.Range.End)
$labelXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:pPr><w:ind w:left="1278"/><w:rPr><w:i/></w:rPr></w:pPr>
<w:r><w:rPr><w:i/></w:rPr><w:t>- Label &amp; Goto</w:t></w:r>
</w:p>
<w:p>
<w:pPr><w:ind w:left="1278"/><w:rPr><w:i/></w:rPr></w:pPr>
<w:r><w:rPr><w:i/></w:rPr><w:t>- Label &amp; Goto in a Diagram</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$labelRange.InsertXML($labelXml)

Write-Output "All edits applied."
